$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# "Clicking routine works. Added dir with the roughsols for the orders
#  not included on the old solution."
#
# Rows 58-61: the Min/Mid/Max-wavelength columns (D/E/F) for these four
# orders previously derived from the fitted-solution formulas; they are
# now pinned to the rough (manually supplied) solution values, and the
# orders are flagged "Auto Done" in column K. Columns G/H/I/J keep their
# existing formulas and simply recompute from the new D/E/F inputs, and
# the upstream C/D/E/F formulas in rows 51-57 ripple through
# automatically because they reference the rows below them.
# -----------------------------------------------------------------------

$ws.Range("D58").Value = 4010.5
$ws.Range("E58").Value = 4077.5
$ws.Range("F58").Value = 4044.7

$ws.Range("D59").Value = 4058.3
$ws.Range("E59").Value = 4126.1
$ws.Range("F59").Value = 4092.8

$ws.Range("D60").Value = 4107.2
$ws.Range("E60").Value = 4175.8
$ws.Range("F60").Value = 4142.2

$ws.Range("D61").Value = 4157.3
$ws.Range("E61").Value = 4226.7
$ws.Range("F61").Value = 4192.7

$ws.Range("K58").Value = "Auto Done"
$ws.Range("K59").Value = "Auto Done"
$ws.Range("K60").Value = "Auto Done"
$ws.Range("K61").Value = "Auto Done"

# -----------------------------------------------------------------------
# Scroll/selection bookkeeping that went along with the edit above (the
# view had scrolled up a little and the active cell ended on the new
# K59 "Auto Done" flag).
# -----------------------------------------------------------------------

$win = $excel.ActiveWindow
$win.ScrollRow = 30
$win.ScrollColumn = 1
$ws.Range("K59").Select()
